# Apply the "marksheet" result update:
#  - fix the summary counters (Right/Wrong/Not-Attempt/Total) for the student
#  - collapse the 3 parallel "Student Ans / Correct Ans" column groups (A/B, D/E, G/H)
#    down to a single group by clearing the now-unused D/E (rows 19-40) and G/H
#    (rows 15-21) cells
#  - populate the "Student Ans" column (A) for every question row with the
#    student's actual answer, colour-coded the same way the summary counters are
#    (green "correctStyle" when it matches the correct answer in B, red
#    "incorrectStyle" when it does not, and left as the neutral "normalStyle"
#    placeholder when the student did not attempt the question)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

# ---- Summary block (rows 10-12) --------------------------------------------
# Row headers (No. / Marking / Total) get the same bordered/centered look as
# the table headers in row 9 ("mtitleStyle").
$ws.Range("A9").Copy()
$ws.Range("A10").PasteSpecial(-4122)
$ws.Range("A11").PasteSpecial(-4122)
$ws.Range("A12").PasteSpecial(-4122)

$ws.Range("B10").Value = 17
$ws.Range("C10").Value = 1
$ws.Range("D10").Value = 10
$ws.Range("E10").Value = 28

$ws.Range("B11").Value = 4
$ws.Range("C11").Value = -1

$ws.Range("B12").Value = 68
$ws.Range("C12").Value = -1
$ws.Range("E12").Value = "67/112"

# ---- Drop the 2nd/3rd "Student Ans / Correct Ans" column groups -----------
$ws.Range("G15:H21").Clear()
$ws.Range("D19:E40").Clear()

# ---- Fill in the student's answer in column A (and D16/D18, which keep
#      their own still-visible "2nd group" Student-Ans cell) for every
#      question ------------------------------------------------------------
$ws.Range("B10").Copy()   # "correctStyle" (green) swatch
$ws.Range("A16").PasteSpecial(-4122)
$ws.Range("D16").PasteSpecial(-4122)
$ws.Range("A18").PasteSpecial(-4122)
$ws.Range("D18").PasteSpecial(-4122)
$ws.Range("A19").PasteSpecial(-4122)
$ws.Range("A21").PasteSpecial(-4122)
$ws.Range("A22").PasteSpecial(-4122)
$ws.Range("A25").PasteSpecial(-4122)
$ws.Range("A27").PasteSpecial(-4122)
$ws.Range("A29").PasteSpecial(-4122)
$ws.Range("A30").PasteSpecial(-4122)
$ws.Range("A31").PasteSpecial(-4122)
$ws.Range("A33").PasteSpecial(-4122)
$ws.Range("A34").PasteSpecial(-4122)
$ws.Range("A37").PasteSpecial(-4122)
$ws.Range("A39").PasteSpecial(-4122)
$ws.Range("A40").PasteSpecial(-4122)

$ws.Range("A16").Value = "Option A"
$ws.Range("D16").Value = "Option A"
$ws.Range("A18").Value = "Option B"
$ws.Range("D18").Value = "Option D"
$ws.Range("A19").Value = "Option C"
$ws.Range("A21").Value = "Option C"
$ws.Range("A22").Value = "Option D"
$ws.Range("A25").Value = "Option A"
$ws.Range("A27").Value = "Option A"
$ws.Range("A29").Value = "Option D"
$ws.Range("A30").Value = "Option B"
$ws.Range("A31").Value = "Option D"
$ws.Range("A33").Value = "Option D"
$ws.Range("A34").Value = "Option B"
$ws.Range("A37").Value = "Option A"
$ws.Range("A39").Value = "Option D"
$ws.Range("A40").Value = "Option D"

$ws.Range("C10").Copy()   # "incorrectStyle" (red) swatch
$ws.Range("A35").PasteSpecial(-4122)
$ws.Range("A35").Value = "Option C"

"done"
